# Actualización automática de tasas-transfi.xlsx

$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note with new Binance rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.36 = 12885.91 pesos`n✅ 12885.91 pesos = 3.34 = 957.87 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- tasas: update rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("O10").Value = 3840
$wsTasas.Range("N12").Value = 3860.9
$wsTasas.Range("O12").Value = 287
